$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.13401358090623
$ws.Range("C2").Value = 8.414878523901196
$ws.Range("D2").Value = 14.12772636425169
$ws.Range("E2").Value = 14.7818993200728
$ws.Range("G2").Value = 3.753245410544127
$ws.Range("J2").Value = 8.891775086349083
$ws.Range("K2").Value = 11.90576790606516
$ws.Range("L2").Value = 12.07889940573057
$ws.Range("M2").Value = 17.255753184908
$ws.Range("O2").Value = 37.39694571420988
$ws.Range("B3").Value = 16.0286072968903
$ws.Range("C3").Value = 8.398002069580572
$ws.Range("D3").Value = 14.13063530701761
$ws.Range("E3").Value = 14.80692786471205
$ws.Range("G3").Value = 3.75550785560903
$ws.Range("J3").Value = 8.897902160824801
$ws.Range("K3").Value = 11.83300948990802
$ws.Range("L3").Value = 12.09362440160219
$ws.Range("M3").Value = 17.25690668050398
$ws.Range("O3").Value = 37.45854251106244
$ws.Range("B4").Value = 15.9670267448289
$ws.Range("C4").Value = 8.387363050817365
$ws.Range("D4").Value = 14.13470455860315
$ws.Range("E4").Value = 14.82383111315714
$ws.Range("G4").Value = 3.756971235900345
$ws.Range("J4").Value = 8.901860621950334
$ws.Range("K4").Value = 11.79055964748418
$ws.Range("L4").Value = 12.10405115917802
$ws.Range("M4").Value = 17.2601062015397
$ws.Range("O4").Value = 37.50188726565717
$ws.Range("B5").Value = 15.94274498683052
$ws.Range("C5").Value = 8.382958288517738
$ws.Range("D5").Value = 14.13693788286986
$ws.Range("E5").Value = 14.83110593873669
$ws.Range("G5").Value = 3.757586300939442
$ws.Range("J5").Value = 8.903523263618293
$ws.Range("K5").Value = 11.77383600641707
$ws.Range("L5").Value = 12.10864906936552
$ws.Range("M5").Value = 17.26203831339098
$ws.Range("O5").Value = 37.52093885029926
$ws.Range("B6").Value = 15.9387627700612
$ws.Range("C6").Value = 8.382222706674222
$ws.Range("D6").Value = 14.13734348154512
$ws.Range("E6").Value = 14.83233728213512
$ws.Range("G6").Value = 3.757689564770538
$ws.Range("J6").Value = 8.903802340218586
$ws.Range("K6").Value = 11.77109423103431
$ws.Range("L6").Value = 12.10943363730904
$ws.Range("M6").Value = 17.26239713537302
$ws.Range("O6").Value = 37.5241861706514
$ws.Range("B7").Value = 15.96669595152762
$ws.Range("C7").Value = 8.387303926709471
$ws.Range("D7").Value = 14.13473234843558
$ws.Range("E7").Value = 14.8239276579935
$ws.Range("G7").Value = 3.756979454983637
$ws.Range("J7").Value = 8.901882844149894
$ws.Range("K7").Value = 11.79033175824974
$ws.Range("L7").Value = 12.10411175477329
$ws.Range("M7").Value = 17.26012971251829
$ws.Range("O7").Value = 37.50213858228934
$ws.Range("B8").Value = 16.09703235453037
$ws.Range("C8").Value = 8.409117168281698
$ws.Range("D8").Value = 14.12825612865699
$ws.Range("E8").Value = 14.79021072950579
$ws.Range("G8").Value = 3.754010131351879
$ws.Range("J8").Value = 8.893847024923508
$ws.Range("K8").Value = 11.88022888931763
$ws.Range("L8").Value = 12.08368932951495
$ws.Range("M8").Value = 17.25563499991999
$ws.Range("O8").Value = 37.4170373144386
$ws.Range("B9").Value = 16.37644412242465
$ws.Range("C9").Value = 8.449683312053882
$ws.Range("D9").Value = 14.13362340805741
$ws.Range("E9").Value = 14.73625558543832
$ws.Range("G9").Value = 3.748773542319833
$ws.Range("J9").Value = 8.879640514691838
$ws.Range("K9").Value = 12.0734340242826
$ws.Range("L9").Value = 12.05461168101697
$ws.Range("M9").Value = 17.26650126389777
$ws.Range("O9").Value = 37.29401833796707
$ws.Range("B10").Value = 16.59472041806252
$ws.Range("C10").Value = 8.478124399778951
$ws.Range("D10").Value = 14.14850573925874
$ws.Range("E10").Value = 14.70400179758109
$ws.Range("G10").Value = 3.745279785430898
$ws.Range("J10").Value = 8.87013948019986
$ws.Range("K10").Value = 12.2246667960539
$ws.Range("L10").Value = 12.03990460480387
$ws.Range("M10").Value = 17.28635667980888
$ws.Range("O10").Value = 37.23041709586622
$ws.Range("B11").Value = 16.6964967128055
$ws.Range("C11").Value = 8.490763792478621
$ws.Range("D11").Value = 14.1576311041609
$ws.Range("E11").Value = 14.69092654679158
$ws.Range("G11").Value = 3.743766344567387
$ws.Range("J11").Value = 8.866018547386094
$ws.Range("K11").Value = 12.29525094950548
$ws.Range("L11").Value = 12.03465193791764
$ws.Range("M11").Value = 17.29793638005568
$ws.Range("O11").Value = 37.20730356586306
$ws.Range("B12").Value = 16.73536359814935
$ws.Range("C12").Value = 8.495506741542584
$ws.Range("D12").Value = 14.16142321269421
$ws.Range("E12").Value = 14.68620443303198
$ws.Range("G12").Value = 3.743204094540521
$ws.Range("J12").Value = 8.864486828885003
$ws.Range("K12").Value = 12.32221618671912
$ws.Range("L12").Value = 12.03286893148851
$ws.Range("M12").Value = 17.3026846078003
$ws.Range("O12").Value = 37.19938787448375
$ws.Range("B13").Value = 16.72697889051903
$ws.Range("C13").Value = 8.494487195006972
$ws.Range("D13").Value = 14.16059158374442
$ws.Range("E13").Value = 14.68721123862645
$ws.Range("G13").Value = 3.743324703142401
$ws.Range("J13").Value = 8.864815433574641
$ws.Range("K13").Value = 12.31639854980432
$ws.Range("L13").Value = 12.03324378011297
$ws.Range("M13").Value = 17.30164588756146
$ws.Range("O13").Value = 37.20105543891387
$ws.Range("B14").Value = 16.6996879488124
$ws.Range("C14").Value = 8.491154867993147
$ws.Range("D14").Value = 14.15793635179841
$ws.Range("E14").Value = 14.69053346488406
$ws.Range("G14").Value = 3.743719870639559
$ws.Range("J14").Value = 8.865891955764086
$ws.Range("K14").Value = 12.29746477362401
$ws.Range("L14").Value = 12.03450112303371
$ws.Range("M14").Value = 17.29831975791318
$ws.Range("O14").Value = 37.20663556439496
$ws.Range("B15").Value = 16.68301305737012
$ws.Range("C15").Value = 8.489108069488704
$ws.Range("D15").Value = 14.15635370389302
$ws.Range("E15").Value = 14.69259825880539
$ws.Range("G15").Value = 3.743963334435874
$ws.Range("J15").Value = 8.866555101959777
$ws.Range("K15").Value = 12.28589748046531
$ws.Range("L15").Value = 12.03529809705299
$ws.Range("M15").Value = 17.29632961832813
$ws.Range("O15").Value = 37.21016254203806
$ws.Range("B16").Value = 16.58811644068956
$ws.Range("C16").Value = 8.477292344942351
$ws.Range("D16").Value = 14.14795659057586
$ws.Range("E16").Value = 14.70488839508759
$ws.Range("G16").Value = 3.745380214653236
$ws.Range("J16").Value = 8.870412828753519
$ws.Range("K16").Value = 12.22008818206087
$ws.Range("L16").Value = 12.04027675198042
$ws.Range("M16").Value = 17.2856509219432
$ws.Range("O16").Value = 37.23204472888371
$ws.Range("B17").Value = 16.53051449463373
$ws.Range("C17").Value = 8.469967033769226
$ws.Range("D17").Value = 14.14340702009668
$ws.Range("E17").Value = 14.7128367368876
$ws.Range("G17").Value = 3.746268821737564
$ws.Range("J17").Value = 8.872830839294739
$ws.Range("K17").Value = 12.1801596845923
$ws.Range("L17").Value = 12.04369875025317
$ws.Range("M17").Value = 17.27975022031083
$ws.Range("O17").Value = 37.24695928901002
$ws.Range("B18").Value = 16.49761878679212
$ws.Range("C18").Value = 8.465725625250347
$ws.Range("D18").Value = 14.14101206459753
$ws.Range("E18").Value = 14.71755877403646
$ws.Range("G18").Value = 3.746787071115059
$ws.Range("J18").Value = 8.874240554742896
$ws.Range("K18").Value = 12.15736337988915
$ws.Range("L18").Value = 12.04580235456199
$ws.Range("M18").Value = 17.27659626044198
$ws.Range("O18").Value = 37.25608546556867
$ws.Range("B19").Value = 16.48652221505734
$ws.Range("C19").Value = 8.464284745352762
$ws.Range("D19").Value = 14.14023933674764
$ws.Range("E19").Value = 14.71918341342428
$ws.Range("G19").Value = 3.746963770487
$ws.Range("J19").Value = 8.874721117183274
$ws.Range("K19").Value = 12.14967467538232
$ws.Range("L19").Value = 12.04653786702656
$ws.Range("M19").Value = 17.27556968827662
$ws.Range("O19").Value = 37.25926949537555
$ws.Range("B20").Value = 16.53662217496217
$ws.Range("C20").Value = 8.470749734894541
$ws.Range("D20").Value = 14.1438683852849
$ws.Range("E20").Value = 14.71197506387397
$ws.Range("G20").Value = 3.746173488804481
$ws.Range("J20").Value = 8.872571478725037
$ws.Range("K20").Value = 12.18439274498462
$ws.Range("L20").Value = 12.04332046801802
$ws.Range("M20").Value = 17.28035354582008
$ws.Range("O20").Value = 37.24531492037895
$ws.Range("B21").Value = 16.70769534572778
$ws.Range("C21").Value = 8.492134831753113
$ws.Range("D21").Value = 14.1587071427568
$ws.Range("E21").Value = 14.68955142989536
$ws.Range("G21").Value = 3.743603506143796
$ws.Range("J21").Value = 8.865574974827915
$ws.Range("K21").Value = 12.30301982782271
$ws.Range("L21").Value = 12.03412622392964
$ws.Range("M21").Value = 17.29928688988254
$ws.Range("O21").Value = 37.20497383297893
$ws.Range("B22").Value = 16.82138974069197
$ws.Range("C22").Value = 8.505858584296444
$ws.Range("D22").Value = 14.17036552914152
$ws.Range("E22").Value = 14.67623201313658
$ws.Range("G22").Value = 3.741987132557923
$ws.Range("J22").Value = 8.861170100033673
$ws.Range("K22").Value = 12.38191855317514
$ws.Range("L22").Value = 12.02931801782379
$ws.Range("M22").Value = 17.3137766316857
$ws.Range("O22").Value = 37.18348661308599
$ws.Range("B23").Value = 16.76054628511764
$ws.Range("C23").Value = 8.498557200223354
$ws.Range("D23").Value = 14.1639646233723
$ws.Range("E23").Value = 14.68321877663482
$ws.Range("G23").Value = 3.742844051526072
$ws.Range("J23").Value = 8.863505759816197
$ws.Range("K23").Value = 12.33969040462531
$ws.Range("L23").Value = 12.03177460893032
$ws.Range("M23").Value = 17.3058506636779
$ws.Range("O23").Value = 37.19450841756176
$ws.Range("B24").Value = 16.53386020379432
$ws.Range("C24").Value = 8.470395969107432
$ws.Range("D24").Value = 14.14365911461987
$ws.Range("E24").Value = 14.71236415140733
$ws.Range("G24").Value = 3.746216565839312
$ws.Range("J24").Value = 8.872688674672569
$ws.Range("K24").Value = 12.18247848117698
$ws.Range("L24").Value = 12.04349106495593
$ws.Range("M24").Value = 17.28008003943208
$ws.Range("O24").Value = 37.24605662120285
$ws.Range("B25").Value = 16.29845963476552
$ws.Range("C25").Value = 8.438947281974063
$ws.Range("D25").Value = 14.1302441161383
$ws.Range("E25").Value = 14.74955254631659
$ws.Range("G25").Value = 3.750127813870235
$ws.Range("J25").Value = 8.88331861475686
$ws.Range("K25").Value = 12.01945934986055
$ws.Range("L25").Value = 12.0613066077257
$ws.Range("M25").Value = 17.26146686978313
$ws.Range("O25").Value = 37.3225973942121
